$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 9894
$ws.Range("D2").Value = 0.05

$ws.Range("C3").Value = 21855
$ws.Range("D3").Value = 0.14

$ws.Range("C4").Value = 5823
$ws.Range("D4").Value = 0.06
